# "agregados pendientes de la semana"
# Adds this week's pending homework/task rows to the "Trabajos pendientes"
# mini-table (columns J:O) on the "Plan de estudio" sheet, and tweaks the
# "Dias para la entrega" count for the existing Bases de datos row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Bases de datos" / "miercoles 27 de mayo"): days-to-deadline 2 -> 1
$ws.Range("L3").Value = 1

# Row 4: Entrepreneur task - "quien es Elon Musk" due Viernes 29 de mayo
$ws.Range("J4").Value = "Entrepreneur"
$ws.Range("K4").Value = "Viernes 29 de mayo"
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = "quien es Elon Musk"
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = "Pendiente"

# Row 5: JS task - "toDoapp" due Lunes 1 de junio
$ws.Range("J5").Value = "JS"
$ws.Range("K5").Value = "Lunes 1 de junio"
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = "toDoapp"
$ws.Range("N5").Value = 4
$ws.Range("O5").Value = "Pendiente"

# Row 6: JS task - "weatherApp" due Lunes 1 de junio
$ws.Range("J6").Value = "JS"
$ws.Range("K6").Value = "Lunes 1 de junio"
$ws.Range("L6").Value = 7
$ws.Range("M6").Value = "weatherApp"
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = "Pendiente"

# The "Estado" column of the newly-filled rows should carry the same
# "Pendiente" highlight formatting already used on row 3 (O3). Copy just
# the format (no values) down onto O4:O6.
[void]$ws.Range("O3").Copy()
[void]$ws.Range("O4:O6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the author's final view/selection state.
[void]$ws.Range("E1").Select()
[void]$ws.Range("N6").Select()
